$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-30 Wednesday" "2025-07-31 Thursday"

Replace-Text "42÷3=14, 0" "13÷7=1, 6"
Replace-Text "41÷4=10, 1" "98÷7=14, 0"
Replace-Text "70÷5=14, 0" "55÷8=6, 7"
Replace-Text "11÷5=2, 1" "36÷3=12, 0"
Replace-Text "98÷4=24, 2" "43÷3=14, 1"

Replace-Text "90÷5=18, 0" "10÷5=2, 0"
Replace-Text "72÷6=12, 0" "46÷8=5, 6"
Replace-Text "20÷8=2, 4" "79÷2=39, 1"
Replace-Text "55÷2=27, 1" "88÷5=17, 3"
Replace-Text "43÷8=5, 3" "73÷3=24, 1"

Replace-Text "15÷8=1, 7" "13÷8=1, 5"
Replace-Text "51÷2=25, 1" "96÷3=32, 0"
Replace-Text "28÷8=3, 4" "83÷8=10, 3"
Replace-Text "78÷9=8, 6" "29÷5=5, 4"
Replace-Text "96÷9=10, 6" "27÷4=6, 3"

Replace-Text "44÷8=5, 4" "45÷7=6, 3"
Replace-Text "80÷7=11, 3" "33÷6=5, 3"
Replace-Text "35÷2=17, 1" "46÷8=5, 6"
Replace-Text "73÷7=10, 3" "19÷4=4, 3"
Replace-Text "19÷3=6, 1" "30÷3=10, 0"

Replace-Text "25÷7=3, 4" "42÷9=4, 6"
Replace-Text "73÷6=12, 1" "63÷3=21, 0"
Replace-Text "17÷8=2, 1" "57÷6=9, 3"
Replace-Text "60÷4=15, 0" "44÷6=7, 2"
Replace-Text "65÷5=13, 0" "43÷9=4, 7"
